$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" date column (C) for rows 2-7 from 2023-09-06 to 2023-09-14
# (Excel serial date 45175 -> 45183)
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
